$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'23.311.84"
$ws.Range('E2').Value = "'  -0.39%  "
$ws.Range('E3').Value = "'  -0.99%  "
$ws.Range('E4').Value = "'  -0.01%  "
$ws.Range('E5').Value = "'  +0.06%  "
$ws.Range('D6').Value = "'302.43"
$ws.Range('E6').Value = "'  -0.88%  "
$ws.Range('D7').Value = "'0.3755"
$ws.Range('E7').Value = "'  +0.90%  "
$ws.Range('B8').Value = "'Cardano"
$ws.Range('C8').Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range('D8').Value = "'0.3614"
$ws.Range('E8').Value = "'  -0.13%  "
$ws.Range('B9').Value = "'OKB"
$ws.Range('C9').Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range('D9').Value = "'51.31"
$ws.Range('E9').Value = "'  -1.33%  "
$ws.Range('D10').Value = "'0.08126"
$ws.Range('E10').Value = "'  +0.20%  "
$ws.Range('D11').Value = "'1.217"
$ws.Range('E11').Value = "'  -2.55%  "
$ws.Range('D12').Value = "'1.001"
$ws.Range('E12').Value = "'  +0.09%  "
$ws.Range('E13').Value = "'  -2.46%  "
$ws.Range('E14').Value = "'  -2.10%  "
$ws.Range('D15').Value = "'0.00001232"
$ws.Range('E15').Value = "'  -2.90%  "
$ws.Range('D16').Value = "'7.254"
$ws.Range('E16').Value = "'  -0.32%  "
$ws.Range('D17').Value = "'1.617.80"
$ws.Range('E17').Value = "'  -0.72%  "
$ws.Range('D18').Value = "'93.94"
$ws.Range('E18').Value = "'  -0.46%  "
$ws.Range('D19').Value = "'0.06929"
$ws.Range('E19').Value = "'  +0.91%  "
$ws.Range('D20').Value = "'17.45"
$ws.Range('E20').Value = "'  -3.54%  "
$ws.Range('D21').Value = "'6.501"
$ws.Range('E21').Value = "'  -0.11%  "
$ws.Range('E22').Value = "'  +0.20%  "
$ws.Range('D23').Value = "'12.47"
$ws.Range('E23').Value = "'  -2.02%  "
$ws.Range('D24').Value = "'23.307.69"
$ws.Range('E24').Value = "'  -0.42%  "
$ws.Range('D25').Value = "'2.475"
$ws.Range('E25').Value = "'  +2.64%  "
$ws.Range('D26').Value = "'3.066"
$ws.Range('E26').Value = "'  +1.80%  "
$ws.Range('D27').Value = "'21.08"
$ws.Range('E27').Value = "'  -0.57%  "
$ws.Range('D28').Value = "'150.21"
$ws.Range('E28').Value = "'  -0.83%  "
$ws.Range('D29').Value = "'5.270"
$ws.Range('E29').Value = "'  +0.01%  "
$ws.Range('D30').Value = "'132.64"
$ws.Range('E30').Value = "'  -2.31%  "
$ws.Range('D31').Value = "'1.797.00"
$ws.Range('E31').Value = "'  -0.71%  "
$ws.Range('D32').Value = "'6.709"
$ws.Range('E32').Value = "'  -0.59%  "
$ws.Range('D33').Value = "'2.125"
$ws.Range('E33').Value = "'  -6.84%  "
$ws.Range('E34').Value = "'  +10.53%  "
$ws.Range('D35').Value = "'11.27"
$ws.Range('E35').Value = "'  +9.43%  "
$ws.Range('D36').Value = "'0.02740"
$ws.Range('E36').Value = "'  -3.24%  "
$ws.Range('D37').Value = "'0.08754"
$ws.Range('E37').Value = "'  -0.24%  "
$ws.Range('D38').Value = "'0.2471"
$ws.Range('E38').Value = "'  -1.91%  "
$ws.Range('D39').Value = "'0.07072"
$ws.Range('E39').Value = "'  -1.71%  "
$ws.Range('D40').Value = "'5.950"
$ws.Range('E40').Value = "'  -1.49%  "
$ws.Range('D41').Value = "'0.6953"
$ws.Range('E41').Value = "'  -1.14%  "
$ws.Range('D42').Value = "'1.326"
$ws.Range('E42').Value = "'  -3.45%  "
$ws.Range('D43').Value = "'15.94"
$ws.Range('E43').Value = "'  -0.59%  "
$ws.Range('E44').Value = "'  -3.60%  "
$ws.Range('D45').Value = "'0.6436"
$ws.Range('E45').Value = "'  -1.01%  "
$ws.Range('E46').Value = "'  +0.11%  "
$ws.Range('E47').Value = "'  -1.27%  "
$ws.Range('E48').Value = "'  -2.76%  "
$ws.Range('D49').Value = "'0.07954"
$ws.Range('E49').Value = "'  -0.13%  "
$ws.Range('D50').Value = "'125.64"
$ws.Range('E50').Value = "'  -2.10%  "
$ws.Range('D51').Value = "'1.179"
$ws.Range('E51').Value = "'  -1.44%  "
